$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '50.948.64'
$ws.Range('E2').Value = '  -16.28%  '
$ws.Range('D3').Value = '2.234.24'
$ws.Range('E3').Value = '  -23.18%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '424.12'
$ws.Range('E5').Value = '  -19.65%  '
$ws.Range('D6').Value = '111.75'
$ws.Range('E6').Value = '  -22.39%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('D8').Value = '0.441'
$ws.Range('E8').Value = '  -19.35%  '
$ws.Range('D9').Value = '2.236.75'
$ws.Range('E9').Value = '  -23.32%  '
$ws.Range('D10').Value = '5.04'
$ws.Range('E10').Value = '  -16.75%  '
$ws.Range('D11').Value = '0.0810'
$ws.Range('E11').Value = '  -24.85%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value = '0.119'
$ws.Range('E12').Value = '  -7.22%  '
$ws.Range('B13').Value = 'Cardano'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D13').Value = '0.287'
$ws.Range('E13').Value = '  -20.03%  '
$ws.Range('D14').Value = '2.635.50'
$ws.Range('E14').Value = '  -22.82%  '
$ws.Range('D15').Value = '51.149.27'
$ws.Range('E15').Value = '  -15.85%  '
$ws.Range('D16').Value = '17.67'
$ws.Range('E16').Value = '  -21.70%  '
$ws.Range('D17').Value = '2.255.71'
$ws.Range('E17').Value = '  -22.73%  '
$ws.Range('D18').Value = '0.0000109'
$ws.Range('E18').Value = '  -22.39%  '
$ws.Range('D19').Value = '3.75'
$ws.Range('E19').Value = '  -24.59%  '
$ws.Range('D20').Value = '282.83'
$ws.Range('E20').Value = '  -20.26%  '
$ws.Range('D21').Value = '0.996'
$ws.Range('E21').Value = '  -0.46%  '
$ws.Range('D22').Value = '5.67'
$ws.Range('E22').Value = '  -1.17%  '
$ws.Range('D23').Value = '8.21'
$ws.Range('E23').Value = '  -29.21%  '
$ws.Range('D24').Value = '4.75'
$ws.Range('E24').Value = '  -27.21%  '
$ws.Range('D25').Value = '0.998'
$ws.Range('D26').Value = '51.53'
$ws.Range('E26').Value = '  -20.81%  '
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = '2.332.63'
$ws.Range('E27').Value = '  -23.07%  '
$ws.Range('B28').Value = 'Polygon'
$ws.Range('C28').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D28').Value = '0.349'
$ws.Range('E28').Value = '  -22.85%  '
$ws.Range('B29').Value = 'USDe'
$ws.Range('C29').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').Value = '0.132'
$ws.Range('E30').Value = '  -25.72%  '
$ws.Range('D31').Value = '6.42'
$ws.Range('E31').Value = '  -18.34%  '
$ws.Range('D32').Value = '142.17'
$ws.Range('E32').Value = '  -7.62%  '
$ws.Range('D33').Value = '0.0₃0592'
$ws.Range('E33').Value = '  -30.98%  '
$ws.Range('D34').Value = '15.91'
$ws.Range('E34').Value = '  -18.93%  '
$ws.Range('E35').Value = '  -24.90%  '
$ws.Range('D36').Value = '4.44'
$ws.Range('E36').Value = '  -20.30%  '
$ws.Range('D38').Value = '0.749'
$ws.Range('E38').Value = '  -24.94%  '
$ws.Range('D39').Value = '3.16'
$ws.Range('E39').Value = '  -27.81%  '
$ws.Range('D40').Value = '31.21'
$ws.Range('E40').Value = '  -16.80%  '
$ws.Range('D41').Value = '0.927'
$ws.Range('E41').Value = '  -22.51%  '
$ws.Range('D43').Value = '0.539'
$ws.Range('E43').Value = '  -17.18%  '
$ws.Range('D44').Value = '0.0471'
$ws.Range('E44').Value = '  -19.14%  '
$ws.Range('D45').Value = '2.94'
$ws.Range('E45').Value = '  -21.24%  '
$ws.Range('D46').Value = '1.817.30'
$ws.Range('E46').Value = '  -20.77%  '
$ws.Range('D47').Value = '1.07'
$ws.Range('E47').Value = '  -27.18%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').Value = '0.0192'
$ws.Range('E48').Value = '  -18.97%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').Value = '0.0771'
$ws.Range('E49').Value = '  -16.02%  '
$ws.Range('E50').Value = '  -5.29%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '14.89'
$ws.Range('E51').Value = '  -26.98%  '
